$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update this run's Price (D) and Volume(1h) (E)
# figures for each listed coin. A handful of the new Price strings are plain
# decimals (e.g. "0.999", "211.79") that Excel would otherwise auto-parse into
# floating-point numbers, so those specific cells are forced to Text first to
# preserve the source string formatting (prices in this sheet are inline
# strings, not real numbers - note the dotted-thousands style like
# "26.220.10" used elsewhere in the same column).

$ws.Range("D2").Value = '26.220.10'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '1.586.43'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.79'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.20'
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("D12").Value = '1.808.93'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '1.567.00'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.91'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '26.213.61'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.87'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.31'
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.60'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.09'
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.18'
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("D33").Value = '1.411.26'
$ws.Range("E33").Value = '  +8.05%  '
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("E37").Value = '  -4.83%  '
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  +5.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.821'
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.952'
$ws.Range("E42").Value = '  -13.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.766'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '1.720.53'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.96'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.38'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("E51").Value = '  -0.35%  '
